$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 185.125
$ws.Range("I33").Value = 181
$ws.Range("K33").Value = 181
$ws.Range("M33").Value = 48

$ws.Range("H40").Value = 7943.4
$ws.Range("I40").Value = 5888.6
$ws.Range("J40").Value = 9998.200000000001
$ws.Range("K40").Value = 5888.6
$ws.Range("L40").Value = 9998.200000000001
$ws.Range("M40").Value = -5713.6
$ws.Range("N40").Value = -10348.2

$ws.Range("H51").Value = 5001
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 5001
$ws.Range("K51").Value = 0
$ws.Range("M51").Value = ""
$ws.Range("N51").Value = -5969

$ws.Range("H68").Value = 74268
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = ""

$ws.Range("H71").Value = 74268
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = ""

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("N76").Value = ""

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("N79").Value = ""

$ws.Range("H80").Value = 274.53845
$ws.Range("I80").Value = 150.625
$ws.Range("J80").Value = 472.8
$ws.Range("K80").Value = 451.875
$ws.Range("L80").Value = 1418.4
$ws.Range("M80").Value = 546.125
$ws.Range("N80").Value = -3414.4

$ws.Range("H83").Value = 274.53845
$ws.Range("I83").Value = 150.625
$ws.Range("J83").Value = 472.8
$ws.Range("K83").Value = 1355.625
$ws.Range("L83").Value = 4255.2
$ws.Range("M83").Value = 3636.375
$ws.Range("N83").Value = -14239.2

$ws.Range("H92").Value = 110.4
$ws.Range("I92").Value = 86.75
$ws.Range("J92").Value = 205
$ws.Range("K92").Value = 86.75
$ws.Range("L92").Value = 205
$ws.Range("M92").Value = 1161.25
$ws.Range("N92").Value = -2701

$ws.Range("H138").Value = 2637.64
$ws.Range("I138").Value = 655.8461
$ws.Range("J138").Value = 4784.5835
$ws.Range("K138").Value = 1967.5383
$ws.Range("L138").Value = 14353.7505
$ws.Range("M138").Value = 3172.4617
$ws.Range("N138").Value = -24633.7505

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2281.1428
$ws.Range("I2").Value = 2528.6667
$ws.Range("J2").Value = 796
$ws.Range("K2").Value = 2528.6667
$ws.Range("L2").Value = 796
$ws.Range("M2").Value = -2415.6667
$ws.Range("N2").Value = -1022

$ws.Range("H88").Value = 1131.125
$ws.Range("I88").Value = 1049.8
$ws.Range("K88").Value = 1049.8
$ws.Range("M88").Value = -643.8

$ws.Range("H91").Value = 1131.125
$ws.Range("I91").Value = 1049.8
$ws.Range("K91").Value = 1049.8
$ws.Range("M91").Value = 354.2

$ws.Range("H97").Value = 1056.2858
$ws.Range("I97").Value = 753
$ws.Range("K97").Value = 753
$ws.Range("M97").Value = -257

$ws.Range("H116").Value = 2281.1428
$ws.Range("I116").Value = 2528.6667
$ws.Range("J116").Value = 796
$ws.Range("K116").Value = 2528.6667
$ws.Range("L116").Value = 796
$ws.Range("M116").Value = -234.6667000000002
$ws.Range("N116").Value = -5384

$ws.Range("H137").Value = 78750
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2281.1428
$ws.Range("I3").Value = 2528.6667
$ws.Range("J3").Value = 796
$ws.Range("K3").Value = 2528.6667
$ws.Range("L3").Value = 796
$ws.Range("M3").Value = -2414.6667
$ws.Range("N3").Value = -1024

$ws.Range("H20").Value = 4699.2
$ws.Range("I20").Value = 4832.3335
$ws.Range("J20").Value = 4499.5
$ws.Range("K20").Value = 4832.3335
$ws.Range("L20").Value = 4499.5
$ws.Range("M20").Value = -4585.3335
$ws.Range("N20").Value = -4993.5

$ws.Range("H86").Value = 2816.6365
$ws.Range("I86").Value = 2098.3
$ws.Range("K86").Value = 2098.3
$ws.Range("M86").Value = -975.3000000000002

$ws.Range("H89").Value = 2816.6365
$ws.Range("I89").Value = 2098.3
$ws.Range("K89").Value = 10491.5
$ws.Range("M89").Value = -4875.5

$ws.Range("H105").Value = 6933625
$ws.Range("I105").Value = 12324833
$ws.Range("J105").Value = 2071.2856
$ws.Range("K105").Value = 12324833
$ws.Range("L105").Value = 2071.2856
$ws.Range("M105").Value = -12323086
$ws.Range("N105").Value = -5565.2856

$ws.Range("H122").Value = 62500
$ws.Range("J122").Value = 62500
$ws.Range("L122").Value = 62500
$ws.Range("N122").Value = -72300

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 23170.9

$ws.Range("H58").Value = 2889.5625
$ws.Range("I58").Value = 2064.8462
$ws.Range("K58").Value = 2064.8462
$ws.Range("M58").Value = -1861.8462

$ws.Range("H60").Value = 3018.4
$ws.Range("I60").Value = 3018.4
$ws.Range("K60").Value = 3018.4
$ws.Range("M60").Value = -2507.4

$ws.Range("H132").Value = 2277
$ws.Range("I132").Value = 2277
$ws.Range("K132").Value = 6831
$ws.Range("M132").Value = -4301

$ws.Range("H136").Value = 2889.5625
$ws.Range("I136").Value = 2064.8462
$ws.Range("K136").Value = 6194.5386
$ws.Range("M136").Value = -3644.5386

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 15010
$ws.Range("I87").Value = 17000
$ws.Range("J87").Value = 12025
$ws.Range("K87").Value = 51000
$ws.Range("L87").Value = 36075
$ws.Range("M87").Value = -49752
$ws.Range("N87").Value = -38571

$ws.Range("H90").Value = 15010
$ws.Range("I90").Value = 17000
$ws.Range("J90").Value = 12025
$ws.Range("K90").Value = 153000
$ws.Range("L90").Value = 108225
$ws.Range("M90").Value = -146760
$ws.Range("N90").Value = -120705

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 39997.875
$ws.Range("J15").Value = 39997.875
$ws.Range("L15").Value = 39997.875
$ws.Range("N15").Value = -40573.875

$ws.Range("H81").Value = 39997.875
$ws.Range("J81").Value = 39997.875
$ws.Range("L81").Value = 39997.875
$ws.Range("N81").Value = -41993.875

$ws.Range("H84").Value = 39997.875
$ws.Range("J84").Value = 39997.875
$ws.Range("L84").Value = 119993.625
$ws.Range("N84").Value = -129977.625

$ws.Range("H97").Value = 1034.2222
$ws.Range("J97").Value = 1461.8
$ws.Range("L97").Value = 1461.8
$ws.Range("N97").Value = -2453.8

$ws.Range("H122").Value = 2041.8889
$ws.Range("I122").Value = 1516.4445
$ws.Range("K122").Value = 4549.333500000001
$ws.Range("M122").Value = -2099.333500000001

$ws.Range("H132").Value = 5022
$ws.Range("I132").Value = 5256.8
$ws.Range("K132").Value = 15770.4
$ws.Range("M132").Value = -13240.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").Value = ""

$ws.Range("H68").Value = 5180.364
$ws.Range("I68").Value = 5372
$ws.Range("J68").Value = 5070.857
$ws.Range("K68").Value = 5372
$ws.Range("L68").Value = 5070.857
$ws.Range("M68").Value = -4623
$ws.Range("N68").Value = -6568.857

$ws.Range("H71").Value = 5180.364
$ws.Range("I71").Value = 5372
$ws.Range("J71").Value = 5070.857
$ws.Range("K71").Value = 26860
$ws.Range("L71").Value = 25354.285
$ws.Range("M71").Value = -23116
$ws.Range("N71").Value = -32842.285

$ws.Range("H122").Value = 2058.1052
$ws.Range("I122").Value = 2067.5715
$ws.Range("J122").Value = 2031.6
$ws.Range("K122").Value = 6202.7145
$ws.Range("L122").Value = 6094.799999999999
$ws.Range("M122").Value = -3752.7145
$ws.Range("N122").Value = -10994.8

$ws.Range("H136").Value = 3470.4285
$ws.Range("I136").Value = 3332.1667
$ws.Range("K136").Value = 9996.500100000001
$ws.Range("M136").Value = -7446.500100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 24386.5
$ws.Range("J80").Value = 24386.5
$ws.Range("L80").Value = 24386.5
$ws.Range("N80").Value = -26382.5

$ws.Range("H83").Value = 24386.5
$ws.Range("J83").Value = 24386.5
$ws.Range("L83").Value = 73159.5
$ws.Range("N83").Value = -83143.5

$ws.Range("H126").Value = 3571
$ws.Range("I126").Value = 2718.8
$ws.Range("K126").Value = 8156.400000000001
$ws.Range("M126").Value = -5686.400000000001

$ws.Range("H132").Value = 2995.6667
$ws.Range("I132").Value = 2993.5
$ws.Range("K132").Value = 8980.5
$ws.Range("M132").Value = -6450.5
